# TC02_Canine_Filter_StudyType-Transcriptomics.xlsx
# Insert a new "StatQuery" column between the existing "query" (A) and
# "dbExcel" (B) columns, shifting dbExcel -> C and WebExcel -> D (and the
# corresponding Neo4jData/WebData filenames in row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; this pushes the old B (dbExcel) -> C and old C
# (WebExcel) -> D, carrying their values/styles/column widths along.
$ws.Columns.Item(2).Insert()

# New header cell for the inserted column.
$ws.Range("B1").Value = "StatQuery"

# Column B (new) should be the same width as column A.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# New query text for the inserted column. Column B already inherited the
# wrap-text style ("Normal 2") from the Insert() shift above, matching A2.
$ws.Range("B2").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_type IN [''Transcriptomics'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# Move selection/top-left view to A2, matching the saved workbook view.
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
